# Insert a new data row at row 30 (pushing existing rows 30-112 down to 31-113)
# and populate it with the new "Madrigal" / "Región Metropolitana" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("30:30").Insert()

$ws.Range("A30").Value = 10
$ws.Range("B30").Value = 'Vega Modelo de Temuco'
$ws.Range("C30").Value = 'La Araucanía'
$ws.Range("D30").Value = (Get-Date -Year 2021 -Month 9 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 100112013
$ws.Range("G30").Value = 'Alcachofa'
$ws.Range("H30").Value = 'Madrigal'
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = 12000
$ws.Range("N30").Value = '$/caja 40 unidades'
$ws.Range("O30").Value = 'Región Metropolitana'
$ws.Range("P30").Value = 300
$ws.Range("Q30").Value = 40
$ws.Range("R30").Value = 'Hortaliza'
